# issue #5: property land done
# Rewrite the "土地" (land) sheet headers to the normalized English schema,
# append the shared metadata columns (property_category .. index), and
# fix a handful of OCR-artifact strings (stray spaces / dashes) across the
# other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. 土地 (land) sheet: rename headers + add metadata columns
# ---------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

$land.Range("B1").Value = "name"
$land.Range("C1").Value = "area"
$land.Range("D1").Value = "share_portion"
$land.Range("E1").Value = "owner"
$land.Range("F1").Value = "register_date"
$land.Range("G1").Value = "register_reason"
$land.Range("H1").Value = "acquire_value"

$land.Range("I1").Value = "property_category"
$land.Range("J1").Value = "category"
$land.Range("K1").Value = "date"
$land.Range("L1").Value = "legislator_name"
$land.Range("M1").Value = "legislator_id"
$land.Range("N1").Value = "source_file"
$land.Range("O1").Value = "index"

# also normalize the stray spaces baked into the scraped text
$land.Range("B2").Value = "臺北市大安區瑞安段二小段08340000地號"
$land.Range("D2").Value = "10000分之202"

$land.Range("I2").Value = "land"
$land.Range("J2").Value = "normal"
$land.Range("K2").Value = "2012-04-19"
$land.Range("L2").Value = "賴士葆"
$land.Range("M2").Value = 866
$land.Range("N2").Value = "tmp9edb1"
$land.Range("O2").Value = 14

$land.Range("B3").Value = "新北市新店區華城二段02140000地號"

$land.Range("I3").Value = "land"
$land.Range("J3").Value = "normal"
$land.Range("K3").Value = "2012-04-19"
$land.Range("L3").Value = "賴士葆"
$land.Range("M3").Value = 866
$land.Range("N3").Value = "tmp9edb1"
$land.Range("O3").Value = 15

# ---------------------------------------------------------------
# 2. 建物 (building) sheet: strip stray spaces/dashes from names
# ---------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")
$building.Range("B2").Value = "臺北市大安區瑞安段二小段02940000建號"
$building.Range("B3").Value = "臺北市大安區瑞安段二小段02983000建號"
$building.Range("H3").Value = "(超過五年地下室停車位）"
$building.Range("B4").Value = "新北市新店區華城二段00053000建號"

# ---------------------------------------------------------------
# 3. 汽車 (car) sheet: strip stray space from date
# ---------------------------------------------------------------
$car = $wb.Worksheets.Item("汽車")
$car.Range("E2").Value = "94年05月26日"

# ---------------------------------------------------------------
# 4. 存款 (deposit) sheet: strip stray spaces from branch names
# ---------------------------------------------------------------
$deposit = $wb.Worksheets.Item("存款")
$deposit.Range("B3").Value = "台北富邦商業銀行和平分行"
$deposit.Range("B5").Value = "合作金庫商業銀行西門支庫"
$deposit.Range("B6").Value = "合作金庫商業銀行東門支庫"
$deposit.Range("B7").Value = "合作金庫商業銀行大安支庫"
$deposit.Range("B8").Value = "台北富邦商業銀行南門分行"
$deposit.Range("B9").Value = "台新國際商業銀行信義分行"
$deposit.Range("B10").Value = "中華郵政股份有限公司台北青田郵局"
$deposit.Range("B12").Value = "台北富邦商業銀行和平分行"

# ---------------------------------------------------------------
# 5. 保險 (insurance) sheet: strip stray space from policy name
# ---------------------------------------------------------------
$insurance = $wb.Worksheets.Item("保險")
$insurance.Range("C2").Value = "安泰人壽靈活理財變額保險甲型"
